$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string must be forced to Text
# format first, otherwise Excel auto-converts them to numbers (losing exact
# text such as trailing zeros, e.g. "43.60" -> 43.6).

$ws.Range("D2").Value = "94.188.95"
$ws.Range("E2").Value = "  +2.11%  "

$ws.Range("D3").Value = "3.075.50"
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.38"
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "607.66"
$ws.Range("E6").Value = "  -1.25%  "

$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.379"
$ws.Range("E8").Value = "  -3.02%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.805"
$ws.Range("E10").Value = "  +8.95%  "

$ws.Range("D11").Value = "3.071.87"
$ws.Range("E11").Value = "  -1.19%  "

$ws.Range("E12").Value = "  -2.31%  "

$ws.Range("D13").Value = "93.810.53"
$ws.Range("E13").Value = "  +1.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000240"
$ws.Range("E14").Value = "  -3.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.69"
$ws.Range("E15").Value = "  -1.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("E16").Value = "  -2.74%  "

$ws.Range("D17").Value = "3.643.74"
$ws.Range("E17").Value = "  -1.50%  "

$ws.Range("D18").Value = "3.067.21"
$ws.Range("E18").Value = "  -1.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.55"
$ws.Range("E19").Value = "  -4.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.30"
$ws.Range("E20").Value = "  -2.50%  "

$ws.Range("E21").Value = "  -0.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "441.46"
$ws.Range("E22").Value = "  -1.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.84"
$ws.Range("E23").Value = "  -5.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000189"
$ws.Range("E24").Value = "  -4.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.37"
$ws.Range("E25").Value = "  +6.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.51"
$ws.Range("E26").Value = "  -4.60%  "

$ws.Range("E27").Value = "  -2.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.87"
$ws.Range("E28").Value = "  +0.34%  "

$ws.Range("D29").Value = "3.231.00"

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("E31").Value = "  +7.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.178"
$ws.Range("E32").Value = "  +4.39%  "

$ws.Range("E33").Value = "  -8.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.87"
$ws.Range("E35").Value = "  -2.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.45"
$ws.Range("E36").Value = "  -5.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("E37").Value = "  -4.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.45"
$ws.Range("E38").Value = "  -2.71%  "

$ws.Range("E39").Value = "  -1.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "483.66"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.82"
$ws.Range("E41").Value = "  -0.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "24.04"
$ws.Range("E42").Value = "  +0.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.434"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.24"
$ws.Range("E44").Value = "  -4.00%  "

$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.07"
$ws.Range("E46").Value = "  -6.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.38"
$ws.Range("E47").Value = "  -1.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.675"
$ws.Range("E48").Value = "  -2.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.82"
$ws.Range("E49").Value = "  -3.92%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.60"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000270"
$ws.Range("E51").Value = "  +9.23%  "
